$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D price updates. Values are stored as plain text strings in the
# sheet (t="inlineStr"), so force the cell format to Text first to avoid
# Excel's automatic number coercion (which would also lose exact
# representations like trailing zeros, e.g. "0.0001500").
$dUpdates = @{
    2  = "275.03"
    3  = "21.18"
    5  = "0.06184"
    6  = "3.575"
    7  = "1.523"
    8  = "6.533"
    9  = "0.8220"
    10 = "0.1651"
    11 = "0.08272"
    12 = "0.03473"
    13 = "0.03164"
    14 = "0.09136"
    15 = "3.762"
    16 = "0.001611"
    17 = "0.04679"
    18 = "0.006264"
    19 = "0.006137"
    20 = "0.001068"
    21 = "0.0001500"
    22 = "3.726"
    23 = "2.321"
    24 = "0.01391"
    25 = "0.3327"
    28 = "0.0002736"
    40 = "0.04733"
    41 = "0.005299"
    42 = "0.007034"
    44 = "0.01159"
    45 = "0.00006056"
    47 = "0.7228"
    49 = "0.00001900"
    50 = "0.01240"
}

foreach ($row in $dUpdates.Keys) {
    $cell = $ws.Range("D$row")
    $cell.NumberFormat = "@"
    $cell.Value = $dUpdates[$row]
}

# Rows 41 and 42 swap places (KickToken <-> CEJI) with updated price/volume data
$ws.Range("B41").Value = "CEJI"
$ws.Range("C41").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("E41").Value = "40CEJICEJIBestin24h"

$ws.Range("B42").Value = "KickToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("E42").Value = "41KickTokenKICK"
